$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first three columns with the new values
$ws.Range("A1").Value = 44.0
$ws.Range("B1").Value = 41.0
$ws.Range("C1").Value = 34.0

$ws.Range("A2").Value = 107.0
$ws.Range("B2").Value = 107.0
$ws.Range("C2").Value = 91.0

$ws.Range("A3").Value = 170.0
$ws.Range("B3").Value = 173.0
$ws.Range("C3").Value = 148.0

# Remove the now-unused fourth column's data
$ws.Range("D1:D3").ClearContents()
